# Revert "Powerpoint writer: consolidate text run nodes."
# Split runs that end with a trailing space into a text-only run and a
# separate single-space run, matching the pre-consolidation OOXML shape.
# We do this by re-assigning Characters(start,len).Text to itself, which
# causes the host to split the underlying <a:r> runs at that boundary
# without introducing any new run-level formatting (rPr stays empty).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape 1: Title "Testing custom properties" ---
$tr1 = $s.Shapes.Item(1).TextFrame.TextRange
# "Testing"(1-7) " "(8) "custom"(9-14) " "(15) "properties"(16-25)
$c = $tr1.Characters(1, 7);  $c.Text = $c.Text
$c = $tr1.Characters(8, 1);  $c.Text = $c.Text
$c = $tr1.Characters(9, 6);  $c.Text = $c.Text
$c = $tr1.Characters(15, 1); $c.Text = $c.Text

# --- Shape 2: Subtitle "This is a subtitle<br/><br/>A. M." ---
$tr2 = $s.Shapes.Item(2).TextFrame.TextRange
# "This"(1-4) " "(5) "is"(6-7) " "(8) "a"(9) " "(10) "subtitle"(11-18)
# <br/>(19) <br/>(20) "A."(21-22) " "(23) "M."(24-25)
$c = $tr2.Characters(1, 4);  $c.Text = $c.Text
$c = $tr2.Characters(5, 1);  $c.Text = $c.Text
$c = $tr2.Characters(6, 2);  $c.Text = $c.Text
$c = $tr2.Characters(8, 1);  $c.Text = $c.Text
$c = $tr2.Characters(9, 1);  $c.Text = $c.Text
$c = $tr2.Characters(10, 1); $c.Text = $c.Text
$c = $tr2.Characters(21, 2); $c.Text = $c.Text
$c = $tr2.Characters(23, 1); $c.Text = $c.Text
